# Update cryptocurrency price/volume data to reflect the latest scrape
# (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price values (column D) are prefixed with a leading
# apostrophe so Excel keeps storing them as text (matching the original
# workbook, where these cells are plain text, not numbers).
$ws.Range("D2").Value = '''59.330.94'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '''2.604.56'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''542.13'
$ws.Range("E5").Value = '  +4.29%  '
$ws.Range("D6").Value = '''141.67'
$ws.Range("E6").Value = '  +1.58%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("E10").Value = '  +2.29%  '
$ws.Range("D11").Value = '''0.336'
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("E12").Value = '  +0.64%  '
$ws.Range("D13").Value = '''3.060.36'
$ws.Range("E13").Value = '  +0.57%  '
$ws.Range("D14").Value = '''59.265.78'
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("D15").Value = '''20.61'
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '''0.0000134'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '''2.582.63'
$ws.Range("E17").Value = '  +0.49%  '
$ws.Range("D18").Value = '''341.42'
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("D19").Value = '''4.35'
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").Value = '''10.14'
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("E21").Value = '  -1.61%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = '''67.72'
$ws.Range("E23").Value = '  +2.29%  '
$ws.Range("E24").Value = '  +1.24%  '
$ws.Range("E25").Value = '  -2.01%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").Value = '''7.23'
$ws.Range("E27").Value = '  +2.87%  '
$ws.Range("D28").Value = '''0.0₃0747'
$ws.Range("E28").Value = '  +3.91%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +7.54%  '
$ws.Range("D31").Value = '''5.81'
$ws.Range("E31").Value = '  -2.23%  '
$ws.Range("D32").Value = '''18.76'
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").Value = '''149.59'
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("D35").Value = '''1.12'
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").Value = '''37.15'
$ws.Range("E36").Value = '  +1.59%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").Value = '''0.835'
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("D39").Value = '''0.827'
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("D40").Value = '''3.56'
$ws.Range("E40").Value = '  +1.78%  '
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("D42").Value = '''275.35'
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  +1.87%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '''0.0956'
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("D47").Value = '''1.951.65'
$ws.Range("E47").Value = '  -1.52%  '
$ws.Range("D48").Value = '''18.56'
$ws.Range("E48").Value = '  +3.82%  '
$ws.Range("E49").Value = '  +1.47%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("E51").Value = '  -0.98%  '
